$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D values are treated as text (not auto-converted to numbers/dates)
$ws.Range('D2:D51').NumberFormat = '@'

$ws.Range('D2').Value = '29.207.46'
$ws.Range('E2').Value = '  -1.28%  '
$ws.Range('D3').Value = '1.858.17'
$ws.Range('E3').Value = '  -0.38%  '
$ws.Range('D4').Value = '0.9996'
$ws.Range('E4').Value = '  -0.39%  '
$ws.Range('B5').Value = 'BNB'
$ws.Range('C5').Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range('D5').Value = '238.37'
$ws.Range('E5').Value = '  -1.32%  '
$ws.Range('B6').Value = 'XRP'
$ws.Range('C6').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range('D6').Value = '0.6915'
$ws.Range('E6').Value = '  -3.04%  '
$ws.Range('D7').Value = '1.000'
$ws.Range('E7').Value = '  -0.39%  '
$ws.Range('D8').Value = '0.07750'
$ws.Range('E8').Value = '  +3.75%  '
$ws.Range('D9').Value = '0.3058'
$ws.Range('E9').Value = '  -2.46%  '
$ws.Range('D10').Value = '23.31'
$ws.Range('E10').Value = '  -3.91%  '
$ws.Range('D11').Value = '0.08070'
$ws.Range('E11').Value = '  -1.18%  '
$ws.Range('D12').Value = '1.856.23'
$ws.Range('E12').Value = '  -0.54%  '
$ws.Range('D13').Value = '0.7237'
$ws.Range('E13').Value = '  -2.37%  '
$ws.Range('D14').Value = '5.209'
$ws.Range('E14').Value = '  -0.97%  '
$ws.Range('D15').Value = '89.55'
$ws.Range('E15').Value = '  -2.24%  '
$ws.Range('D16').Value = '29.218.28'
$ws.Range('E16').Value = '  -1.96%  '
$ws.Range('D17').Value = '5.753'
$ws.Range('E17').Value = '  -3.46%  '
$ws.Range('D18').Value = '0.000007820'
$ws.Range('E18').Value = '  -0.46%  '
$ws.Range('E19').Value = '  -0.83%  '
$ws.Range('D20').Value = '235.50'
$ws.Range('E20').Value = '  -3.63%  '
$ws.Range('D21').Value = '0.9998'
$ws.Range('E21').Value = '  -0.20%  '
$ws.Range('D22').Value = '2.106.50'
$ws.Range('E22').Value = '  -1.43%  '
$ws.Range('D23').Value = '1.0000'
$ws.Range('E23').Value = '  -0.24%  '
$ws.Range('E24').Value = '  -2.75%  '
$ws.Range('D25').Value = '162.03'
$ws.Range('E25').Value = '  -0.87%  '
$ws.Range('D26').Value = '8.985'
$ws.Range('E26').Value = '  -1.68%  '
$ws.Range('D27').Value = '0.1443'
$ws.Range('E27').Value = '  -2.92%  '
$ws.Range('D28').Value = '18.10'
$ws.Range('E28').Value = '  -1.92%  '
$ws.Range('D29').Value = '1.962'
$ws.Range('E29').Value = '  -1.47%  '
$ws.Range('D30').Value = '1.403'
$ws.Range('E30').Value = '  -1.28%  '
$ws.Range('D31').Value = '4.521'
$ws.Range('E31').Value = '  +0.69%  '
$ws.Range('D32').Value = '1.490'
$ws.Range('E32').Value = '  -2.13%  '
$ws.Range('D33').Value = '4.026'
$ws.Range('E33').Value = '  -2.69%  '
$ws.Range('D34').Value = '0.05187'
$ws.Range('E34').Value = '  -3.91%  '
$ws.Range('D35').Value = '1.187'
$ws.Range('E35').Value = '  -2.72%  '
$ws.Range('D36').Value = '0.7058'
$ws.Range('E36').Value = '  -3.82%  '
$ws.Range('D37').Value = '1.027'
$ws.Range('E37').Value = '  +2.98%  '
$ws.Range('D38').Value = '2.671'
$ws.Range('E38').Value = '  -1.08%  '
$ws.Range('D39').Value = '0.01852'
$ws.Range('E39').Value = '  -2.74%  '
$ws.Range('D40').Value = '2.682'
$ws.Range('E40').Value = '  -1.56%  '
$ws.Range('D41').Value = '0.9235'
$ws.Range('E41').Value = '  +4.98%  '
$ws.Range('D42').Value = '1.096.55'
$ws.Range('E42').Value = '  +6.04%  '
$ws.Range('D43').Value = '5.977'
$ws.Range('E43').Value = '  +0.09%  '
$ws.Range('D44').Value = '0.4293'
$ws.Range('E44').Value = '  -2.87%  '
$ws.Range('D45').Value = '70.58'
$ws.Range('E45').Value = '  -0.75%  '
$ws.Range('D46').Value = '0.9998'
$ws.Range('E46').Value = '  -0.34%  '
$ws.Range('D47').Value = '102.38'
$ws.Range('E47').Value = '  -0.68%  '
$ws.Range('D48').Value = '1.793'
$ws.Range('D49').Value = '2.003.74'
$ws.Range('E49').Value = '  -1.47%  '
$ws.Range('D50').Value = '9.191'
$ws.Range('E50').Value = '  -2.38%  '
$ws.Range('D51').Value = '7.012'
$ws.Range('E51').Value = '  -5.00%  '

# Restore original (default) style now that text values are locked in
$ws.Range('D2:D51').Style = 'Normal'
